$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.749791622161865
$ws.Range("B1").Value = 1.992416262626648
$ws.Range("C1").Value = 2.156231164932251
$ws.Range("D1").Value = 2.843078851699829
$ws.Range("E1").Value = 1.736220955848694
